$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily profit row (row 88) below the existing data
# (previously A1:B87, now A1:B88).
#
# Column A stores the date as literal text (matching every other row in
# the sheet, e.g. A87 = "11/12/2025"), so we briefly force a text number
# format before assigning the value - otherwise Excel's COM layer
# auto-parses "11/13/2025" into a date serial number. ClearFormats()
# afterwards drops the cell back to the default (unstyled) format, same
# as the rest of the sheet, without re-triggering the date parser.
$ws.Range("A88").NumberFormat = "@"
$ws.Range("A88").Value = "11/13/2025"
$ws.Range("A88").ClearFormats()

$ws.Range("B88").Value = 10080.44
